$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/3/2025  Through  11/9/2025"

# --- Data table updates (rows 14-33) ---
# Numeric cells
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("N14").Value = -84.615384615384
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = -37.5
$ws.Range("I15").Value = 44
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 10
$ws.Range("L15").Value = 37.5
$ws.Range("M15").Value = 69.230769230769
$ws.Range("N15").Value = -20
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 45
$ws.Range("G16").Value = 45
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 445
$ws.Range("J16").Value = 463
$ws.Range("K16").Value = -3.887688984881
$ws.Range("L16").Value = 10.972568578553
$ws.Range("M16").Value = 2.534562211981
$ws.Range("N16").Value = -69.520547945205
$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = 58.333333333333
$ws.Range("F17").Value = 76
$ws.Range("G17").Value = 45
$ws.Range("H17").Value = 68.888888888888
$ws.Range("I17").Value = 809
$ws.Range("J17").Value = 683
$ws.Range("K17").Value = 18.448023426061
$ws.Range("L17").Value = 36.886632825719
$ws.Range("M17").Value = 109.585492227979
$ws.Range("N17").Value = 45.503597122302
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 100
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = 37.5
$ws.Range("I18").Value = 199
$ws.Range("J18").Value = 176
$ws.Range("K18").Value = 13.068181818181
$ws.Range("L18").Value = 1.530612244897
$ws.Range("M18").Value = -42.651296829971
$ws.Range("N18").Value = -90.124069478908
$ws.Range("C19").Value = 21
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = -12.5
$ws.Range("F19").Value = 75
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -2.597402597402
$ws.Range("I19").Value = 795
$ws.Range("J19").Value = 765
$ws.Range("K19").Value = 3.92156862745
$ws.Range("L19").Value = 26.996805111821
$ws.Range("M19").Value = 65.625
$ws.Range("N19").Value = 0.505689001264
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 17
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = -41.379310344827
$ws.Range("I20").Value = 279
$ws.Range("J20").Value = 306
$ws.Range("K20").Value = -8.823529411764
$ws.Range("L20").Value = -5.423728813559
$ws.Range("M20").Value = 78.846153846153
$ws.Range("N20").Value = -78.488820354664
$ws.Range("C21").Value = 59
$ws.Range("D21").Value = 55
$ws.Range("E21").Value = 7.272727272727
$ws.Range("F21").Value = 241
$ws.Range("G21").Value = 221
$ws.Range("H21").Value = 9.049773755656
$ws.Range("I21").Value = 2577
$ws.Range("J21").Value = 2444
$ws.Range("K21").Value = 5.441898527004
$ws.Range("L21").Value = 19.860465116279
$ws.Range("M21").Value = 40.28307022319
$ws.Range("N21").Value = -58.522452921294
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = 40
$ws.Range("I22").Value = 42
$ws.Range("J22").Value = 42
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -6.666666666666
$ws.Range("M22").Value = 75
$ws.Range("C24").Value = 64
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = 45.454545454545
$ws.Range("F24").Value = 166
$ws.Range("G24").Value = 146
$ws.Range("H24").Value = 13.698630136986
$ws.Range("I24").Value = 1478
$ws.Range("J24").Value = 1607
$ws.Range("K24").Value = -8.027380211574
$ws.Range("L24").Value = -33.543165467625
$ws.Range("M24").Value = 35.347985347985
$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 47.826086956521
$ws.Range("F25").Value = 94
$ws.Range("G25").Value = 79
$ws.Range("H25").Value = 18.987341772151
$ws.Range("I25").Value = 633
$ws.Range("J25").Value = 867
$ws.Range("K25").Value = -26.989619377162
$ws.Range("L25").Value = -56.643835616438
$ws.Range("C26").Value = 20
$ws.Range("D26").Value = 16
$ws.Range("E26").Value = 25
$ws.Range("F26").Value = 86
$ws.Range("G26").Value = 80
$ws.Range("H26").Value = 7.5
$ws.Range("I26").Value = 812
$ws.Range("J26").Value = 777
$ws.Range("K26").Value = 4.504504504504
$ws.Range("L26").Value = 14.366197183098
$ws.Range("M26").Value = 4.639175257731
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = -66.666666666666
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 56
$ws.Range("J27").Value = 59
$ws.Range("K27").Value = -5.084745762711
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 300
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 400
$ws.Range("I28").Value = 116
$ws.Range("J28").Value = 93
$ws.Range("K28").Value = 24.731182795698
$ws.Range("L28").Value = 31.818181818181
$ws.Range("C29").Value = 2
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 16
$ws.Range("J29").Value = 36
$ws.Range("K29").Value = -55.555555555555
$ws.Range("L29").Value = 6.666666666666
$ws.Range("M29").Value = -52.941176470588
$ws.Range("N29").Value = -80.246913580246
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 12
$ws.Range("J30").Value = 21
$ws.Range("K30").Value = -42.857142857142
$ws.Range("M30").Value = -53.846153846153
$ws.Range("N30").Value = -84
$ws.Range("L31").Value = -66.666666666666

# Text "0" placeholder cells (copy style+shared-string from stable source D30)
$ws.Range("D30").Copy($ws.Range("C14"))
$ws.Range("D30").Copy($ws.Range("F33"))
